# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions update).
# For each changed cell we set the literal text exactly as scraped from coinranking.com:
#  - "Price" (col D) and "Volume(1h)" (col E) are always text, not numbers, because the
#    site formats them with locale dot-grouping ("67.248.31") or padded percent strings
#    ("  -8.43%  "). Cells whose new value would otherwise auto-parse as a plain number
#    (e.g. "566.31") are pinned to the Text number format first so Excel keeps them as
#    strings instead of silently converting them to floats.
#  - Rows 28/29, 40/41 and 48/49 additionally swapped their Coin/Link/Price/Volume values
#    (ranking reorder), so B/C are rewritten too for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.248.31'
$ws.Range('E2').Value = '  -8.43%  '

$ws.Range('D3').Value = '3.677.44'
$ws.Range('E3').Value = '  -7.65%  '

$ws.Range('E4').Value = '  +0.34%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.31'
$ws.Range('E5').Value = '  -7.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.11'
$ws.Range('E6').Value = '  +0.75%  '

$ws.Range('D7').Value = '3.674.15'
$ws.Range('E7').Value = '  -7.52%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.617'
$ws.Range('E8').Value = '  -9.70%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.25%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.696'
$ws.Range('E10').Value = '  -11.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.161'
$ws.Range('E11').Value = '  -14.21%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.34'
$ws.Range('E12').Value = '  -8.72%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000289'
$ws.Range('E13').Value = '  -14.66%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.40'
$ws.Range('E14').Value = '  -7.76%  '

$ws.Range('D15').Value = '4.287.28'
$ws.Range('E15').Value = '  -7.10%  '

$ws.Range('D16').Value = '3.696.29'
$ws.Range('E16').Value = '  -7.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.126'
$ws.Range('E17').Value = '  -3.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.07'
$ws.Range('E18').Value = '  -8.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.11'
$ws.Range('E19').Value = '  -9.95%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.66'
$ws.Range('E20').Value = '  -11.52%  '

$ws.Range('D21').Value = '67.405.65'
$ws.Range('E21').Value = '  -8.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '400.80'
$ws.Range('E22').Value = '  -11.79%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.44'
$ws.Range('E23').Value = '  -8.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.90'
$ws.Range('E24').Value = '  -9.68%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.00'
$ws.Range('E25').Value = '  -11.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.53'
$ws.Range('E26').Value = '  -12.05%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.46'
$ws.Range('E27').Value = '  -5.74%  '

$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.96'
$ws.Range('E28').Value = '  -0.35%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.73'
$ws.Range('E29').Value = '  -11.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.34'
$ws.Range('E30').Value = '  -11.63%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.23'
$ws.Range('E31').Value = '  -11.55%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.43'
$ws.Range('E32').Value = '  -7.10%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.35'
$ws.Range('E33').Value = '  -11.45%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.115'
$ws.Range('E34').Value = '  -11.27%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.35'
$ws.Range('E35').Value = '  -9.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.55'
$ws.Range('E36').Value = '  -11.55%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '587.74'
$ws.Range('E37').Value = '  -9.41%  '

$ws.Range('D38').Value = '0.0₃0865'
$ws.Range('E38').Value = '  -18.54%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.16%  '

$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.389'
$ws.Range('E41').Value = '  -9.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.133'
$ws.Range('E42').Value = '  -9.56%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.98'
$ws.Range('E43').Value = '  -12.37%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0432'
$ws.Range('E44').Value = '  -10.58%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.81'
$ws.Range('E45').Value = '  -14.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.51'
$ws.Range('E46').Value = '  -3.32%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.08'
$ws.Range('E47').Value = '  -14.85%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.133'
$ws.Range('E48').Value = '  -11.14%  '

$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.68'
$ws.Range('E49').Value = '  -14.68%  '

$ws.Range('D50').Value = '2.705.53'
$ws.Range('E50').Value = '  -3.74%  '

$ws.Range('E51').Value = '  -10.27%  '
